$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the cryptocurrency price/volume table (rows 2-51) to the
# latest scraped values. Row 38/39 also had their coin identity
# (name/link) swapped along with the numeric data.
# NumberFormat is forced to text ("@") before each write so that
# price strings like "607.61" or "3.131.40" are preserved verbatim
# instead of being auto-converted to numbers by Excel.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.322.97"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.63%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.658.10"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.38%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "607.61"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.08%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "152.21"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +5.36%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.85%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +1.93%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +6.36%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.31%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.74%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "28.09"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +2.61%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.131.40"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.37%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "64.176.00"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.66%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000149"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +2.40%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.651.81"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.84%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.18"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +6.72%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.65"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +4.62%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "347.24"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.13%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.95"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.39%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.09%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.58"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.03%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "66.73"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.61%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.76"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +14.12%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.72"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +4.69%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +8.75%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +4.15%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "556.35"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +1.70%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.164"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.13%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.02%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +1.75%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0₃0867"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +6.90%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.78"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.76%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.36"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +5.09%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "169.05"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -2.03%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.409"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.67%  "
$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = "Stacks"
$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.97"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +7.01%  "
$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = "FirstDigitalUSD"
$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.998"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.23%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "19.40"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +1.39%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.05%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "167.21"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -2.67%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "40.33"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.52%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +2.68%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.05%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "22.07"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.38%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.36%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.02"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +15.46%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0247"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +2.95%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.38%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +1.23%  "
